$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '52.277.01'
$ws.Range('E2').Value = '  +1.05%  '

# Row 3
$ws.Range('D3').Value = '2.810.59'
$ws.Range('E3').Value = '  +1.85%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.71%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '116.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.57%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.550'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.54%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.36%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.45%  '

# Row 11
$ws.Range('E11').Value = '  +3.63%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.00'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.36%  '

# Row 13
$ws.Range('E13').Value = '  +1.42%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.23%  '

# Row 15
$ws.Range('D15').Value = '3.251.21'
$ws.Range('E15').Value = '  +2.05%  '

# Row 16
$ws.Range('D16').Value = '2.803.02'
$ws.Range('E16').Value = '  +1.79%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.899'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.05%  '

# Row 18
$ws.Range('D18').Value = '52.268.65'
$ws.Range('E18').Value = '  +1.18%  '

# Row 19
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.59%  '

# Row 20
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.76%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.16%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0985'
$ws.Range('E22').Value = '  +1.81%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.15%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.98%  '

# Row 25
$ws.Range('E25').Value = '  +3.72%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.07%  '

# Row 27
$ws.Range('E27').Value = '  +0.02%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.61%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.96%  '

# Row 30
$ws.Range('E30').Value = '  -0.22%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.00%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0462'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +33.37%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.21%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.52%  '

# Row 35
$ws.Range('E35').Value = '  +1.37%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.31%  '

# Row 37
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.01%  '

# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.97'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.96%  '

# Row 39
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.70%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.47%  '

# Row 41
$ws.Range('E41').Value = '  +10.39%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.77%  '

# Row 43
$ws.Range('E43').Value = '  +1.88%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.01%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '

# Row 46
$ws.Range('E46').Value = '  -1.20%  '

# Row 47
$ws.Range('D47').Value = '2.059.60'
$ws.Range('E47').Value = '  -2.74%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.48%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.960'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.63%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.18%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.05%  '
